$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.289615
$ws.Range("H2").Value = 27.868845
$ws.Range("I2").Value = 0.3997772888571298
$ws.Range("J2").Value = 0.439676128114975
$ws.Range("M2").Value = 9.363528666666667
$ws.Range("N2").Value = 28.090586
$ws.Range("O2").Value = 0.04175743631338733
$ws.Range("P2").Value = 0.04324026421082073
$ws.Range("Q2").Value = 86.98357635479667
$ws.Range("R2").Value = 782.85218719317
$ws.Range("S2").Value = 0.01669367467899025
$ws.Range("T2").Value = 0.01901171194688218

$ws.Range("G3").Value = 9.289615
$ws.Range("H3").Value = 27.868845
$ws.Range("I3").Value = 0.3997772888571298
$ws.Range("J3").Value = 0.439676128114975
$ws.Range("O3").Value = 0.1749266505387075
$ws.Range("P3").Value = 0.1811383852696593
$ws.Range("Q3").Value = 364.3840955519667
$ws.Range("R3").Value = 3279.4568599677
$ws.Range("S3").Value = 0.06993170210122304
$ws.Range("T3").Value = 0.07964222388836241

$ws.Range("G4").Value = 9.289615
$ws.Range("H4").Value = 27.868845
$ws.Range("I4").Value = 0.3997772888571298
$ws.Range("J4").Value = 0.439676128114975
$ws.Range("M4").Value = 71.284935
$ws.Range("N4").Value = 213.854805
$ws.Range("O4").Value = 0.3179011075133629
$ws.Range("P4").Value = 0.3291899382573772
$ws.Range("Q4").Value = 662.209601450025
$ws.Range("R4").Value = 5959.886413050225
$ws.Range("S4").Value = 0.1270896428863711
$ws.Range("T4").Value = 0.1447369574674113

$ws.Range("G5").Value = 9.289615
$ws.Range("H5").Value = 27.868845
$ws.Range("I5").Value = 0.3997772888571298
$ws.Range("J5").Value = 0.439676128114975
$ws.Range("M5").Value = 23.0690325
$ws.Range("N5").Value = 46.138065
$ws.Range("O5").Value = 0.1028782726814826
$ws.Range("P5").Value = 0.07102102180339065
$ws.Range("Q5").Value = 214.3024303474875
$ws.Range("R5").Value = 1285.814582084925
$ws.Range("S5").Value = 0.04112839693490763
$ws.Range("T5").Value = 0.03122624788128402

$ws.Range("G6").Value = 9.289615
$ws.Range("H6").Value = 27.868845
$ws.Range("I6").Value = 0.3997772888571298
$ws.Range("J6").Value = 0.439676128114975
$ws.Range("M6").Value = 81.293813
$ws.Range("N6").Value = 243.881439
$ws.Range("O6").Value = 0.3625365329530597
$ws.Range("P6").Value = 0.3754103904587522
$ws.Range("Q6").Value = 755.1882246519949
$ws.Range("R6").Value = 6796.694021867955
$ws.Range("S6").Value = 0.1449338722556377
$ws.Range("T6").Value = 0.1650589869310351

$ws.Range("G7").Value = 5.424084000000001
$ws.Range("I7").Value = 0.2334247001682348
$ws.Range("J7").Value = 0.2567211075690851
$ws.Range("M7").Value = 9.363528666666667
$ws.Range("N7").Value = 28.090586
$ws.Range("O7").Value = 0.04175743631338733
$ws.Range("P7").Value = 0.04324026421082073
$ws.Range("Q7").Value = 50.78856602440801
$ws.Range("R7").Value = 457.0970942196721
$ws.Range("S7").Value = 0.009747217051246596
$ws.Range("T7").Value = 0.01110068851978177

$ws.Range("G8").Value = 5.424084000000001
$ws.Range("I8").Value = 0.2334247001682348
$ws.Range("J8").Value = 0.2567211075690851
$ws.Range("O8").Value = 0.1749266505387075
$ws.Range("P8").Value = 0.1811383852696593
$ws.Range("S8").Value = 0.04083220095343137
$ws.Range("T8").Value = 0.04650204688970257

$ws.Range("G9").Value = 5.424084000000001
$ws.Range("I9").Value = 0.2334247001682348
$ws.Range("J9").Value = 0.2567211075690851
$ws.Range("M9").Value = 71.284935
$ws.Range("N9").Value = 213.854805
$ws.Range("O9").Value = 0.3179011075133629
$ws.Range("P9").Value = 0.3291899382573772
$ws.Range("Q9").Value = 386.6554753745401
$ws.Range("R9").Value = 3479.89927837086
$ws.Range("S9").Value = 0.0742059707044565
$ws.Range("T9").Value = 0.0845100055500326

$ws.Range("G10").Value = 5.424084000000001
$ws.Range("I10").Value = 0.2334247001682348
$ws.Range("J10").Value = 0.2567211075690851
$ws.Range("M10").Value = 23.0690325
$ws.Range("N10").Value = 46.138065
$ws.Range("O10").Value = 0.1028782726814826
$ws.Range("P10").Value = 0.07102102180339065
$ws.Range("Q10").Value = 125.12837007873
$ws.Range("R10").Value = 750.7702204723801
$ws.Range("S10").Value = 0.02401432995450098
$ws.Range("T10").Value = 0.01823259537805458

$ws.Range("G11").Value = 5.424084000000001
$ws.Range("I11").Value = 0.2334247001682348
$ws.Range("J11").Value = 0.2567211075690851
$ws.Range("M11").Value = 81.293813
$ws.Range("N11").Value = 243.881439
$ws.Range("O11").Value = 0.3625365329530597
$ws.Range("P11").Value = 0.3754103904587522
$ws.Range("Q11").Value = 440.944470392292
$ws.Range("R11").Value = 3968.500233530629
$ws.Range("S11").Value = 0.08462498150459932
$ws.Range("T11").Value = 0.09637577123151354

$ws.Range("G12").Value = 1.032076666666667
$ws.Range("H12").Value = 3.09623
$ws.Range("I12").Value = 0.04441527573453837
$ws.Range("J12").Value = 0.04884803866659809
$ws.Range("M12").Value = 9.363528666666667
$ws.Range("N12").Value = 28.090586
$ws.Range("O12").Value = 0.04175743631338733
$ws.Range("P12").Value = 0.04324026421082073
$ws.Range("Q12").Value = 9.663879454531113
$ws.Range("R12").Value = 86.97491509078002
$ws.Range("S12").Value = 0.001854668047826524
$ws.Range("T12").Value = 0.002112202098124089

$ws.Range("G13").Value = 1.032076666666667
$ws.Range("H13").Value = 3.09623
$ws.Range("I13").Value = 0.04441527573453837
$ws.Range("J13").Value = 0.04884803866659809
$ws.Range("O13").Value = 0.1749266505387075
$ws.Range("P13").Value = 0.1811383852696593
$ws.Range("Q13").Value = 40.48309028131112
$ws.Range("R13").Value = 364.3478125318001
$ws.Range("S13").Value = 0.007769415416995927
$ws.Range("T13").Value = 0.00884825484765746

$ws.Range("G14").Value = 1.032076666666667
$ws.Range("H14").Value = 3.09623
$ws.Range("I14").Value = 0.04441527573453837
$ws.Range("J14").Value = 0.04884803866659809
$ws.Range("M14").Value = 71.284935
$ws.Range("N14").Value = 213.854805
$ws.Range("O14").Value = 0.3179011075133629
$ws.Range("P14").Value = 0.3291899382573772
$ws.Range("Q14").Value = 73.57151809835001
$ws.Range("R14").Value = 662.1436628851501
$ws.Range("S14").Value = 0.01411966534652114
$ws.Range("T14").Value = 0.0160802828326514

$ws.Range("G15").Value = 1.032076666666667
$ws.Range("H15").Value = 3.09623
$ws.Range("I15").Value = 0.04441527573453837
$ws.Range("J15").Value = 0.04884803866659809
$ws.Range("M15").Value = 23.0690325
$ws.Range("N15").Value = 46.138065
$ws.Range("O15").Value = 0.1028782726814826
$ws.Range("P15").Value = 0.07102102180339065
$ws.Range("Q15").Value = 23.809010165825
$ws.Range("R15").Value = 142.85406099495
$ws.Range("S15").Value = 0.004569366848241076
$ws.Range("T15").Value = 0.003469237619193332

$ws.Range("G16").Value = 1.032076666666667
$ws.Range("H16").Value = 3.09623
$ws.Range("I16").Value = 0.04441527573453837
$ws.Range("J16").Value = 0.04884803866659809
$ws.Range("M16").Value = 81.293813
$ws.Range("N16").Value = 243.881439
$ws.Range("O16").Value = 0.3625365329530597
$ws.Range("P16").Value = 0.3754103904587522
$ws.Range("Q16").Value = 83.90144754166334
$ws.Range("R16").Value = 755.11302787497
$ws.Range("S16").Value = 0.0161021600749537
$ws.Range("T16").Value = 0.01833806126897181

$ws.Range("G17").Value = 6.325986
$ws.Range("H17").Value = 12.651972
$ws.Range("I17").Value = 0.2722379272368294
$ws.Range("J17").Value = 0.1996053321183234
$ws.Range("M17").Value = 9.363528666666667
$ws.Range("N17").Value = 28.090586
$ws.Range("O17").Value = 0.04175743631338733
$ws.Range("P17").Value = 0.04324026421082073
$ws.Range("Q17").Value = 59.23355125593201
$ws.Range("R17").Value = 355.401307535592
$ws.Range("S17").Value = 0.01136795790868048
$ws.Range("T17").Value = 0.008630987298684924

$ws.Range("G18").Value = 6.325986
$ws.Range("H18").Value = 12.651972
$ws.Range("I18").Value = 0.2722379272368294
$ws.Range("J18").Value = 0.1996053321183234
$ws.Range("O18").Value = 0.1749266505387075
$ws.Range("P18").Value = 0.1811383852696593
$ws.Range("Q18").Value = 248.13608390492
$ws.Range("R18").Value = 1488.81650342952
$ws.Range("S18").Value = 0.04762166876113893
$ws.Range("T18").Value = 0.03615618755112716

$ws.Range("G19").Value = 6.325986
$ws.Range("H19").Value = 12.651972
$ws.Range("I19").Value = 0.2722379272368294
$ws.Range("J19").Value = 0.1996053321183234
$ws.Range("M19").Value = 71.284935
$ws.Range("N19").Value = 213.854805
$ws.Range("O19").Value = 0.3179011075133629
$ws.Range("P19").Value = 0.3291899382573772
$ws.Range("Q19").Value = 450.9475008209101
$ws.Range("R19").Value = 2705.68500492546
$ws.Range("S19").Value = 0.08654473857573039
$ws.Range("T19").Value = 0.06570806695587414

$ws.Range("G20").Value = 6.325986
$ws.Range("H20").Value = 12.651972
$ws.Range("I20").Value = 0.2722379272368294
$ws.Range("J20").Value = 0.1996053321183234
$ws.Range("M20").Value = 23.0690325
$ws.Range("N20").Value = 46.138065
$ws.Range("O20").Value = 0.1028782726814826
$ws.Range("P20").Value = 0.07102102180339065
$ws.Range("Q20").Value = 145.934376628545
$ws.Range("R20").Value = 583.73750651418
$ws.Range("S20").Value = 0.02800736771251216
$ws.Range("T20").Value = 0.01417617464444848

$ws.Range("G21").Value = 6.325986
$ws.Range("H21").Value = 12.651972
$ws.Range("I21").Value = 0.2722379272368294
$ws.Range("J21").Value = 0.1996053321183234
$ws.Range("M21").Value = 81.293813
$ws.Range("N21").Value = 243.881439
$ws.Range("O21").Value = 0.3625365329530597
$ws.Range("P21").Value = 0.3754103904587522
$ws.Range("Q21").Value = 514.263522924618
$ws.Range("R21").Value = 3085.581137547708
$ws.Range("S21").Value = 0.09869619427876748
$ws.Range("T21").Value = 0.07493391566818869

$ws.Range("G22").Value = 1.165213666666667
$ws.Range("H22").Value = 3.495641
$ws.Range("I22").Value = 0.05014480800326766
$ws.Range("J22").Value = 0.05514939353101857
$ws.Range("M22").Value = 9.363528666666667
$ws.Range("N22").Value = 28.090586
$ws.Range("O22").Value = 0.04175743631338733
$ws.Range("P22").Value = 0.04324026421082073
$ws.Range("Q22").Value = 10.91051157062511
$ws.Range("R22").Value = 98.19460413562601
$ws.Range("S22").Value = 0.002093918626643485
$ws.Range("T22").Value = 0.002384674347347771

$ws.Range("G23").Value = 1.165213666666667
$ws.Range("H23").Value = 3.495641
$ws.Range("I23").Value = 0.05014480800326766
$ws.Range("J23").Value = 0.05514939353101857
$ws.Range("O23").Value = 0.1749266505387075
$ws.Range("P23").Value = 0.1811383852696593
$ws.Range("Q23").Value = 45.70537401745111
$ws.Range("R23").Value = 411.34836615706
$ws.Range("S23").Value = 0.008771663305918184
$ws.Range("T23").Value = 0.009989672092809697

$ws.Range("G24").Value = 1.165213666666667
$ws.Range("H24").Value = 3.495641
$ws.Range("I24").Value = 0.05014480800326766
$ws.Range("J24").Value = 0.05514939353101857
$ws.Range("M24").Value = 71.284935
$ws.Range("N24").Value = 213.854805
$ws.Range("O24").Value = 0.3179011075133629
$ws.Range("P24").Value = 0.3291899382573772
$ws.Range("Q24").Value = 83.062180489445
$ws.Range("R24").Value = 747.559624405005
$ws.Range("S24").Value = 0.01594109000028374
$ws.Range("T24").Value = 0.0181546254514078

$ws.Range("G25").Value = 1.165213666666667
$ws.Range("H25").Value = 3.495641
$ws.Range("I25").Value = 0.05014480800326766
$ws.Range("J25").Value = 0.05514939353101857
$ws.Range("M25").Value = 23.0690325
$ws.Range("N25").Value = 46.138065
$ws.Range("O25").Value = 0.1028782726814826
$ws.Range("P25").Value = 0.07102102180339065
$ws.Range("Q25").Value = 26.8803519457775
$ws.Range("R25").Value = 161.282111674665
$ws.Range("S25").Value = 0.005158811231320761
$ws.Range("T25").Value = 0.003916766280410241

$ws.Range("G26").Value = 1.165213666666667
$ws.Range("H26").Value = 3.495641
$ws.Range("I26").Value = 0.05014480800326766
$ws.Range("J26").Value = 0.05514939353101857
$ws.Range("M26").Value = 81.293813
$ws.Range("N26").Value = 243.881439
$ws.Range("O26").Value = 0.3625365329530597
$ws.Range("P26").Value = 0.3754103904587522
$ws.Range("Q26").Value = 94.72466192304432
$ws.Range("R26").Value = 852.521957307399
$ws.Range("S26").Value = 0.0181793248391015
$ws.Range("T26").Value = 0.02070365535904306
